$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "状态编号"
for ($i = 3; $i -le 29; $i++) {
    $ws.Cells.Item($i, 1).Value = $i - 2
}

$ws.Range("I9").Select()

$ws.PageSetup.PaperSize = 9          # xlPaperA4
$ws.PageSetup.Orientation = $xlPortrait
